$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 10
$ws.Range("E3").Value = 10
$ws.Range("E4").Value = 10
$ws.Range("E5").Value = 10

$ws.Range("E6").Select()

$wb.Save()
